# Applies the district-name corrections to "the official names from website"
# and removes a handful of stray empty Address cells, as described in the
# commit message / diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- District name (column G) corrections ---------------------------------
$districtFixes = @{
    "G4"  = "Davangere"
    "G5"  = "Kalaburagi (Gulbarga)"
    "G7"  = "Kalaburagi (Gulbarga)"
    "G9"  = "Vijayapura (Bijapur)"
    "G11" = "Vijayapura (Bijapur)"
    "G14" = "Chikkamagaluru (Chikmagalur)"
    "G24" = "Dharwad"
    "G25" = "Kalaburagi (Gulbarga)"
    "G27" = "Kalaburagi (Gulbarga)"
    "G36" = "Davangere"
    "G43" = "Kalaburagi (Gulbarga)"
    "G49" = "Kalaburagi (Gulbarga)"
    "G51" = "Davangere"
    "G54" = "Vijayapura (Bijapur)"
    "G55" = "Vijayapura (Bijapur)"
    "G56" = "Kalaburagi (Gulbarga)"
    "G57" = "Kalaburagi (Gulbarga)"
}

foreach ($addr in $districtFixes.Keys) {
    $ws.Range($addr).Value = $districtFixes[$addr]
}

# --- Remove stray empty Address cells (column F) ---------------------------
$emptyAddressCells = @("F6", "F13", "F20", "F21", "F28", "F30", "F33", "F37", "F48", "F50")

foreach ($addr in $emptyAddressCells) {
    $ws.Range($addr).ClearContents()
}
